$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1428571428571428
$ws.Range("C2").Value = 0.7142857142857143
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("J3").Value = 0.4
$ws.Range("P3").Value = 0.4
$ws.Range("S3").Value = 0.2
$ws.Range("P4").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("D6").Value = 0.1666666666666667
$ws.Range("F6").Value = 0.1666666666666667
$ws.Range("Q6").Value = 0.6666666666666666
$ws.Range("B7").Value = 0.1666666666666667
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.3333333333333333
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.0625
$ws.Range("F8").Value = 0.0625
$ws.Range("Q8").Value = 0.3125
$ws.Range("R8").Value = 0.1875
$ws.Range("S8").Value = 0.375
$ws.Range("B9").Value = 0.1428571428571428
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.5714285714285714
$ws.Range("B10").Value = 0.03703703703703703
$ws.Range("E10").Value = 0.01851851851851852
$ws.Range("F10").Value = 0.05555555555555555
$ws.Range("J10").Value = 0.09259259259259259
$ws.Range("Q10").Value = 0.2962962962962963
$ws.Range("R10").Value = 0.1296296296296296
$ws.Range("S10").Value = 0.3703703703703703
$ws.Range("G11").Value = 0.2
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.3
$ws.Range("L11").Value = 0.4
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.25
$ws.Range("G13").Value = 1
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.1666666666666667
$ws.Range("J15").Value = 0.6666666666666666
$ws.Range("H16").Value = 0.2
$ws.Range("J16").Value = 0.8
$ws.Range("H17").Value = 0.16
$ws.Range("I17").Value = 0.16
$ws.Range("J17").Value = 0.44
$ws.Range("K17").Value = 0.12
$ws.Range("O17").Value = 0.04
$ws.Range("S17").Value = 0.08
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("J18").Value = 0.5833333333333334
$ws.Range("K18").Value = 0.25
$ws.Range("F19").Value = 0.025
$ws.Range("H19").Value = 0.2
$ws.Range("I19").Value = 0.05
$ws.Range("J19").Value = 0.425
$ws.Range("K19").Value = 0.05
$ws.Range("M19").Value = 0.025
$ws.Range("O19").Value = 0.075
$ws.Range("S19").Value = 0.15
